# main figures draft 1
# Adds a new observation row (row 14, "Rainer desktop" / 9 strains / 58 days)
# to the "data" sheet, mirroring the pattern of rows 2-13, and updates the
# worksheet selection to reflect where the user ended up after entering it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- New data row (row 14) -------------------------------------------------
# Columns A-C / E / J are plain input values; D stays at the same
# wait-time value used by the other "Rainer desktop" rows; the rest are
# formulas copied down from row 13, re-pointed at row 14.

$ws.Range("A14").Value2 = "Rainer desktop"
$ws.Range("B14").Value2 = 9
$ws.Range("C14").Value2 = 58
$ws.Range("D14").Value2 = 1000000
$ws.Range("E14").Value2 = 100

$ws.Range("F14").Formula = "=D14*E14"
$ws.Range("G14").Formula = "= IF(C14<12,C14, 12)"
$ws.Range("H14").Formula = "=C14*D14*E14"
$ws.Range("I14").Formula = "=H14/G14"

$ws.Range("J14").Value2 = 55000

$ws.Range("K14").Formula = "=I14*regression!B`$18+regression!B`$17"
$ws.Range("L14").Formula = "=K14/60"
$ws.Range("M14").Formula = "=L14/60"
$ws.Range("N14").Formula = "=M14/24"
$ws.Range("N14").Style = "Normal"
$ws.Range("O14").Formula = "=I14*regression!B`$18+regression!B`$17"

$ws.Range("S14").Formula = "=LOG10(I14)"
$ws.Range("T14").Formula = "=LOG10(J14)"

# --- Selection / view -------------------------------------------------------
$ws.Activate()
$ws.Range("T16").Select()
